$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.508.31"
$ws.Range("E2").Value = "'  -4.38%  "
$ws.Range("D3").Value = "'3.064.98"
$ws.Range("E3").Value = "'  -3.85%  "
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'535.49"
$ws.Range("E5").Value = "'  -5.75%  "
$ws.Range("D6").Value = "'132.88"
$ws.Range("E6").Value = "'  -10.19%  "
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("D8").Value = "'3.057.72"
$ws.Range("E8").Value = "'  -3.97%  "
$ws.Range("E9").Value = "'  -3.61%  "
$ws.Range("E10").Value = "'  -4.04%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "'  -11.61%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "'  -5.14%  "
$ws.Range("E13").Value = "'  -2.59%  "
$ws.Range("D14").Value = "'34.20"
$ws.Range("E14").Value = "'  -9.99%  "
$ws.Range("D15").Value = "'3.527.21"
$ws.Range("E15").Value = "'  -4.68%  "
$ws.Range("D16").Value = "'62.655.38"
$ws.Range("E16").Value = "'  -4.32%  "
$ws.Range("E17").Value = "'  -3.07%  "
$ws.Range("D18").Value = "'3.072.10"
$ws.Range("E18").Value = "'  -3.76%  "
$ws.Range("D19").Value = "'6.55"
$ws.Range("E19").Value = "'  -6.70%  "
$ws.Range("D20").Value = "'475.78"
$ws.Range("E20").Value = "'  -10.07%  "
$ws.Range("D21").Value = "'13.22"
$ws.Range("E21").Value = "'  -7.92%  "
$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "'  -5.22%  "
$ws.Range("D23").Value = "'7.10"
$ws.Range("E23").Value = "'  -7.35%  "
$ws.Range("D24").Value = "'78.42"
$ws.Range("E24").Value = "'  -2.26%  "
$ws.Range("D25").Value = "'11.95"
$ws.Range("E25").Value = "'  -9.36%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  -0.02%  "
$ws.Range("D27").Value = "'2.68"
$ws.Range("E27").Value = "'  -7.29%  "
$ws.Range("D28").Value = "'8.20"
$ws.Range("E28").Value = "'  -11.54%  "
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("D30").Value = "'25.60"
$ws.Range("E30").Value = "'  -5.51%  "
$ws.Range("D31").Value = "'1.85"
$ws.Range("E31").Value = "'  -16.94%  "
$ws.Range("D32").Value = "'1.09"
$ws.Range("E32").Value = "'  -5.98%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.38"
$ws.Range("E33").Value = "'  -10.98%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'56.79"
$ws.Range("E34").Value = "'  +3.97%  "
$ws.Range("D35").Value = "'5.90"
$ws.Range("E35").Value = "'  -5.59%  "
$ws.Range("D36").Value = "'5.18"
$ws.Range("E36").Value = "'  -6.44%  "
$ws.Range("E37").Value = "'  -14.88%  "
$ws.Range("D38").Value = "'3.103.56"
$ws.Range("E38").Value = "'  -2.96%  "
$ws.Range("E39").Value = "'  -12.20%  "
$ws.Range("D40").Value = "'0.0785"
$ws.Range("E40").Value = "'  -7.00%  "
$ws.Range("D41").Value = "'7.99"
$ws.Range("E41").Value = "'  -5.91%  "
$ws.Range("D42").Value = "'0.112"
$ws.Range("E42").Value = "'  -11.24%  "
$ws.Range("E43").Value = "'  -9.95%  "
$ws.Range("E44").Value = "'  -0.03%  "
$ws.Range("D45").Value = "'0.248"
$ws.Range("E45").Value = "'  -10.19%  "
$ws.Range("E46").Value = "'  -12.12%  "
$ws.Range("D47").Value = "'24.25"
$ws.Range("E47").Value = "'  -7.35%  "
$ws.Range("D48").Value = "'117.81"
$ws.Range("E48").Value = "'  -4.49%  "
$ws.Range("E49").Value = "'  -3.87%  "
$ws.Range("D50").Value = "'0.0₃0508"
$ws.Range("D51").Value = "'1.98"
$ws.Range("E51").Value = "'  -8.63%  "
